$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge "THU Oct 11" + " 13:42:06 IST 2018" into a single run.
#    Re-assigning the paragraph's text (via a temporary different
#    value, then the final value) collapses the two existing runs
#    into a single run while keeping the Courier New run formatting
#    that was already present.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$dateIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t.StartsWith("THU Oct 11")) {
        $dateIdx = $i
        break
    }
}

if ($dateIdx -gt 0) {
    $r = $d.Paragraphs.Item($dateIdx).Range
    $r.MoveEnd(1, -1)
    $r.Text = "THU Oct 11 13:42:06 IST 2018#TMP#"

    $r2 = $d.Paragraphs.Item($dateIdx).Range
    $r2.MoveEnd(1, -1)
    $r2.Text = "THU Oct 11 13:42:06 IST 2018"
}

# ---------------------------------------------------------------------
# 2) Locate the LAST "Amount Received mode ... - CASH AND CLEARD"
#    paragraph (the 3rd / final occurrence in the document) and append
#    the new purchase-record block right after it.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = $paras.Count; $i -ge 1; $i--) {
    $t = $paras.Item($i).Range.Text
    if ($t -match "CASH AND CLEARD") {
        $targetIndex = $i
        break
    }
}

$anchorEnd = $d.Paragraphs.Item($targetIndex).Range.End
$insertionPoint = $d.Range($anchorEnd, $anchorEnd)

$MARK = [char]1

$lines = @(
    $MARK,
    "SAT Oct 13 14:00:11 IST 2018",
    ("Person Name" + "`t`t`t`t" + "- HANUMANTHARAYA"),
    ("Bill number" + "`t`t`t`t" + "- 8385"),
    "---------------------------------------------------------------",
    ("Item Name" + "`t`t`t`t" + "- CARROT"),
    ("Number of Pockets" + "`t`t`t" + "- 1"),
    ("Number of KGs" + "`t`t`t" + "- 91"),
    ("Rate" + "`t`t`t`t`t" + "- 30"),
    ("Total Price" + "`t`t`t`t" + "- 2730.0"),
    ("Amount balance" + "`t`t`t" + "- 2730.0"),
    $MARK,
    ("Item Name" + "`t`t`t`t" + "- CARROT"),
    ("Amount Received" + "`t`t`t" + "- 2500"),
    ("Amount balance" + "`t`t`t" + "- 230.0"),
    ("Amount Received mode" + "`t`t" + "- CASH"),
    $MARK,
    $MARK
)

$blockText = [string]::Join("`r", $lines)
$insertionPoint.InsertAfter($blockText)

# ---------------------------------------------------------------------
# 3) Re-resolve the freshly created paragraphs, clear the placeholder
#    marker character from the paragraphs that must stay empty, and
#    apply the red / bold run formatting that two of the new lines
#    need ("Amount Received" -> red, bold "Amount balance" -> bold).
# ---------------------------------------------------------------------
$startIndex = $targetIndex + 1

for ($k = 0; $k -lt $lines.Count; $k++) {
    if ($lines[$k] -eq $MARK) {
        $idx = $startIndex + $k
        $r = $d.Paragraphs.Item($idx).Range
        $r.MoveEnd(1, -1)
        $r.Text = ""
    }
}

# "Amount Received" (red) is lines[13]; "Amount balance" (bold) is lines[14]
$redIdx = $startIndex + 13
$redRange = $d.Paragraphs.Item($redIdx).Range
$redRange.MoveEnd(1, -1)
$redRange.Font.Color = 255

$boldIdx = $startIndex + 14
$boldRange = $d.Paragraphs.Item($boldIdx).Range
$boldRange.MoveEnd(1, -1)
$boldRange.Font.Bold = 1

Write-Host "Done. dateIdx=$dateIdx targetIndex=$targetIndex startIndex=$startIndex totalParas=$($d.Paragraphs.Count)"
